$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H4").Value = 748.52
$ws1.Range("H19").Value = "1 de 17"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 748.52
$ws2.Range("F19").Value = 22150.19

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D7").Value = 748.52
$ws3.Range("E7").Value = 51.48000000000002
$ws3.Range("F7").Value = 0.93565
$ws3.Range("D19").Value = 22150.19
$ws3.Range("E19").Value = 25069.11386304603
$ws3.Range("F19").Value = 0.4690918371910773
